$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max"/1) entirely; this shifts D->C and E->D
$ws.Range("C:C").Delete() | Out-Null

# Update B2 with the new value
$ws.Range("B2").Value = 119770.3808655355
